$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 7: G8 program (shortened high-school track) ---
$ws.Range("A7").Value = "G8"
$ws.Range("B7").Value = "G8"
$ws.Range("C7").Value = 2011
$ws.Range("E7").Value = 18
$ws.Rows.Item(7).RowHeight = 45
$ws.Range("F7").WrapText = $true

# Unify the "Information" / "Mentoring" categories into "Education"
$ws.Range("D5").Value = "Education"
$ws.Range("D6").Value = "Education"
$ws.Range("D7").Value = "Education"

# --- New row 8: Tuition fees program ---
$ws.Range("A8").Value = "tuitionFees"
$ws.Range("B8").Value = "Tuition Fees"
$ws.Range("C8").Value = 2008
$ws.Range("D8").Value = "Education"
$ws.Range("E8").Value = 21
$ws.Rows.Item(8).RowHeight = 75
$ws.Range("F8").WrapText = $true
$ws.Range("F8").Value = "In 2006 and 2007 some German federal states introduced tuition fees of about 1000€ per year. All tuition fees in Germany have since been abolished again. Lower Saxony was the last state to scrap tuition fees in 2015."

# G8 description entered last
$ws.Range("F7").Value = 'Between 2007 and 2013 some German federal states shortened the length of the academic high school track "Gymnasium" from 9 to 8 years.'

$ws.Range("G2").Select()
